$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 675
$ws.Range("J38").Value = 2150
$ws.Range("L38").Value = 6450
$ws.Range("N38").Value = -7194
$ws.Range("H113").Value = 2985.7144
$ws.Range("J113").Value = 3380
$ws.Range("L113").Value = 3380
$ws.Range("N113").Value = -9888
$ws.Range("H116").Value = 18184746
$ws.Range("I116").Value = 40002200
$ws.Range("J116").Value = 3534.1667
$ws.Range("K116").Value = 40002200
$ws.Range("L116").Value = 3534.1667
$ws.Range("M116").Value = -39998758
$ws.Range("N116").Value = -10418.1667
$ws.Range("H126").Value = 44855.5
$ws.Range("J126").Value = 44855.5
$ws.Range("L126").Value = 44855.5
$ws.Range("N126").Value = -54735.5
$ws.Range("H132").Value = 2166.4465
$ws.Range("I132").Value = 1726.6097
$ws.Range("J132").Value = 3368.6667
$ws.Range("K132").Value = 5179.8291
$ws.Range("L132").Value = 10106.0001
$ws.Range("M132").Value = -2649.8291
$ws.Range("N132").Value = -15166.0001
$ws.Range("H138").Value = 5058.846
$ws.Range("I138").Value = 1655.5454
$ws.Range("J138").Value = 9463.117
$ws.Range("K138").Value = 4966.6362
$ws.Range("L138").Value = 28389.351
$ws.Range("M138").Value = 173.3638000000001
$ws.Range("N138").Value = -38669.351

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 50000
$ws.Range("I3").Value = 50000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 50000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -49885
$ws.Range("N3").Value = ""
$ws.Range("H32").Value = 9294.281000000001
$ws.Range("I32").Value = 10941.895
$ws.Range("K32").Value = 10941.895
$ws.Range("M32").Value = -10654.895
$ws.Range("H45").Value = 1177.1428
$ws.Range("I45").Value = 948
$ws.Range("J45").Value = 1750
$ws.Range("K45").Value = 948
$ws.Range("L45").Value = 1750
$ws.Range("M45").Value = -571
$ws.Range("N45").Value = -2504
$ws.Range("H61").Value = 1525.1072
$ws.Range("I61").Value = 1373.1052
$ws.Range("J61").Value = 1846
$ws.Range("K61").Value = 1373.1052
$ws.Range("L61").Value = 1846
$ws.Range("M61").Value = -1161.1052
$ws.Range("N61").Value = -2270
$ws.Range("H110").Value = 957.17145
$ws.Range("I110").Value = 1012.85187
$ws.Range("J110").Value = 769.25
$ws.Range("K110").Value = 1012.85187
$ws.Range("L110").Value = 769.25
$ws.Range("M110").Value = 1032.14813
$ws.Range("N110").Value = -4859.25
$ws.Range("H132").Value = 2325.795
$ws.Range("I132").Value = 1765.4286
$ws.Range("J132").Value = 2639.6
$ws.Range("K132").Value = 5296.2858
$ws.Range("L132").Value = 7918.799999999999
$ws.Range("M132").Value = -2766.2858
$ws.Range("N132").Value = -12978.8
$ws.Range("H136").Value = 1525.1072
$ws.Range("I136").Value = 1373.1052
$ws.Range("J136").Value = 1846
$ws.Range("K136").Value = 4119.3156
$ws.Range("L136").Value = 5538
$ws.Range("M136").Value = -1569.3156
$ws.Range("N136").Value = -10638
$ws.Range("H129").Value = 49926
$ws.Range("J129").Value = 49926
$ws.Range("L129").Value = 49926
$ws.Range("N129").Value = -59926

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 39847.332
$ws.Range("J106").Value = 39847.332
$ws.Range("L106").Value = 39847.332
$ws.Range("N106").Value = -42371.332
$ws.Range("H107").Value = 31838.176
$ws.Range("I107").Value = 40680.08
$ws.Range("J107").Value = 3102
$ws.Range("K107").Value = 40680.08
$ws.Range("L107").Value = 3102
$ws.Range("M107").Value = -38760.08
$ws.Range("N107").Value = -6942
$ws.Range("H134").Value = 2262.162
$ws.Range("I134").Value = 2012.7858
$ws.Range("J134").Value = 3038
$ws.Range("K134").Value = 6038.357400000001
$ws.Range("L134").Value = 9114
$ws.Range("M134").Value = -3503.357400000001
$ws.Range("N134").Value = -14184
$ws.Range("H13").Value = 10000
$ws.Range("J13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("N13").Value = -10280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1695.6666
$ws.Range("I16").Value = 1615.25
$ws.Range("J16").Value = 1760
$ws.Range("K16").Value = 1615.25
$ws.Range("L16").Value = 1760
$ws.Range("M16").Value = -1328.25
$ws.Range("N16").Value = -2334
$ws.Range("H31").Value = 2177.923
$ws.Range("I31").Value = 1427.1111
$ws.Range("J31").Value = 3867.25
$ws.Range("K31").Value = 1427.1111
$ws.Range("L31").Value = 3867.25
$ws.Range("M31").Value = -1132.1111
$ws.Range("N31").Value = -4457.25
$ws.Range("H34").Value = 2177.923
$ws.Range("I34").Value = 1427.1111
$ws.Range("J34").Value = 3867.25
$ws.Range("K34").Value = 1427.1111
$ws.Range("L34").Value = 3867.25
$ws.Range("M34").Value = -1225.1111
$ws.Range("N34").Value = -4271.25
$ws.Range("H99").Value = 5678.5713
$ws.Range("I99").Value = 7750
$ws.Range("J99").Value = 2916.6667
$ws.Range("K99").Value = 7750
$ws.Range("L99").Value = 2916.6667
$ws.Range("M99").Value = -6252
$ws.Range("N99").Value = -5912.6667
$ws.Range("H113").Value = 1695.6666
$ws.Range("I113").Value = 1615.25
$ws.Range("J113").Value = 1760
$ws.Range("K113").Value = 1615.25
$ws.Range("L113").Value = 1760
$ws.Range("M113").Value = 554.75
$ws.Range("N113").Value = -6100
$ws.Range("H126").Value = 5678.5713
$ws.Range("I126").Value = 7750
$ws.Range("J126").Value = 2916.6667
$ws.Range("K126").Value = 23250
$ws.Range("L126").Value = 8750.000100000001
$ws.Range("M126").Value = -20780
$ws.Range("N126").Value = -13690.0001
$ws.Range("H132").Value = 376949.28
$ws.Range("I132").Value = 615644.3
$ws.Range("J132").Value = 1857.0714
$ws.Range("K132").Value = 1846932.9
$ws.Range("L132").Value = 5571.2142
$ws.Range("M132").Value = -1844402.9
$ws.Range("N132").Value = -10631.2142
$ws.Range("H134").Value = 1665.8055
$ws.Range("I134").Value = 1440.6923
$ws.Range("K134").Value = 4322.0769
$ws.Range("M134").Value = -1787.0769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 657.2308
$ws.Range("I107").Value = 548.8889
$ws.Range("J107").Value = 901
$ws.Range("K107").Value = 1646.6667
$ws.Range("L107").Value = 2703
$ws.Range("M107").Value = 273.3332999999998
$ws.Range("N107").Value = -6543
$ws.Range("H131").Value = 17874920
$ws.Range("J131").Value = 20853820
$ws.Range("L131").Value = 62561460
$ws.Range("N131").Value = -62571540

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1007131.56
$ws.Range("I3").Value = 1336311
$ws.Range("J3").Value = 19593.334
$ws.Range("K3").Value = 1336311
$ws.Range("L3").Value = 19593.334
$ws.Range("M3").Value = -1336195
$ws.Range("N3").Value = -19825.334
$ws.Range("H11").Value = 42891630
$ws.Range("I11").Value = 51460000
$ws.Range("J11").Value = 49800
$ws.Range("K11").Value = 51460000
$ws.Range("L11").Value = 49800
$ws.Range("M11").Value = -51459861
$ws.Range("N11").Value = -50078
$ws.Range("H111").Value = 40597.668
$ws.Range("J111").Value = 40597.668
$ws.Range("L111").Value = 40597.668
$ws.Range("N111").Value = -46731.668
$ws.Range("H113").Value = 1933.2354
$ws.Range("J113").Value = 2322.0833
$ws.Range("L113").Value = 2322.0833
$ws.Range("N113").Value = -6662.0833
$ws.Range("H126").Value = 1503.7273
$ws.Range("I126").Value = 1171.2222
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 3513.6666
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -1043.6666
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 1926.6316
$ws.Range("I132").Value = 1117.4546
$ws.Range("J132").Value = 3039.25
$ws.Range("K132").Value = 3352.3638
$ws.Range("L132").Value = 9117.75
$ws.Range("M132").Value = -822.3638000000001
$ws.Range("N132").Value = -14177.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 12506536
$ws.Range("I122").Value = 11910634
$ws.Range("K122").Value = 35731902
$ws.Range("M122").Value = -35729452
$ws.Range("H132").Value = 3325.347
$ws.Range("I132").Value = 2866.353
$ws.Range("J132").Value = 4365.7334
$ws.Range("K132").Value = 8599.059000000001
$ws.Range("L132").Value = 13097.2002
$ws.Range("M132").Value = -6069.059000000001
$ws.Range("N132").Value = -18157.2002
$ws.Range("H136").Value = 2611.5952
$ws.Range("I136").Value = 2761.7188
$ws.Range("K136").Value = 8285.1564
$ws.Range("M136").Value = -5735.1564

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 6580
$ws.Range("I29").Value = 6580
$ws.Range("K29").Value = 6580
$ws.Range("M29").Value = -6290
$ws.Range("H46").Value = 66307.94
$ws.Range("J46").Value = 66307.94
$ws.Range("L46").Value = 66307.94
$ws.Range("N46").Value = -66769.94
$ws.Range("H81").Value = 63755.125
$ws.Range("I81").Value = 63755.125
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 127510.25
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -126449.25
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 63755.125
$ws.Range("I84").Value = 63755.125
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 637551.25
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -632247.25
$ws.Range("N84").Value = ""
$ws.Range("H107").Value = 736.4375
$ws.Range("I107").Value = 676.63635
$ws.Range("J107").Value = 868
$ws.Range("K107").Value = 2029.90905
$ws.Range("L107").Value = 2604
$ws.Range("M107").Value = -109.90905
$ws.Range("N107").Value = -6444
$ws.Range("H109").Value = 34125.668
$ws.Range("J109").Value = 34125.668
$ws.Range("L109").Value = 34125.668
$ws.Range("N109").Value = -36899.668
$ws.Range("H122").Value = 13022873
$ws.Range("I122").Value = 16668405
$ws.Range("J122").Value = 6946986
$ws.Range("K122").Value = 50005215
$ws.Range("L122").Value = 20840958
$ws.Range("M122").Value = -50002765
$ws.Range("N122").Value = -20845858
$ws.Range("H132").Value = 973.54095
$ws.Range("I132").Value = 648.38776
$ws.Range("K132").Value = 1945.16328
$ws.Range("M132").Value = 584.8367200000002
$ws.Range("H134").Value = 66307.94
$ws.Range("J134").Value = 66307.94
$ws.Range("L134").Value = 198923.82
$ws.Range("N134").Value = -203993.82
$ws.Range("H136").Value = 2285.1562
$ws.Range("I136").Value = 1934.375
$ws.Range("J136").Value = 3337.5
$ws.Range("K136").Value = 5803.125
$ws.Range("L136").Value = 10012.5
$ws.Range("M136").Value = -3253.125
$ws.Range("N136").Value = -15112.5

Write-Output "Applied all Asura_Profits updates"